$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) and Volume(1h) (E) columns to Text format so that
# numeric-looking strings (e.g. "1.512", "0.08850") are preserved as
# literal text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "24.745.58"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.701.00"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "316.12"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.3935"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "0.4049"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "1.512"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").Value = "1.002"
$ws.Range("E10").Value = "  +0.13%  "
$ws.Range("D11").Value = "52.64"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "0.08850"
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").Value = "7.623"
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("D14").Value = "23.66"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "8.109"
$ws.Range("E15").Value = "  +7.07%  "
$ws.Range("D16").Value = "0.00001323"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Value = "1.692.63"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "99.42"
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").Value = "0.07081"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "19.83"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "7.112"
$ws.Range("E21").Value = "  +3.76%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "14.75"
$ws.Range("E23").Value = "  +5.07%  "
$ws.Range("D24").Value = "24.730.22"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").Value = "3.139"
$ws.Range("E25").Value = "  +3.10%  "
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").Value = "22.76"
$ws.Range("E27").Value = "  +2.13%  "
$ws.Range("D28").Value = "9.192"
$ws.Range("E28").Value = "  +22.29%  "
$ws.Range("D29").Value = "164.54"
$ws.Range("E29").Value = "  +2.20%  "
$ws.Range("D30").Value = "136.14"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("D31").Value = "5.141"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("D32").Value = "7.971"
$ws.Range("E32").Value = "  +9.76%  "
$ws.Range("D33").Value = "0.09044"
$ws.Range("E33").Value = "  +6.03%  "
$ws.Range("D34").Value = "1.073"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "0.03001"
$ws.Range("E35").Value = "  +9.36%  "
$ws.Range("D36").Value = "0.2782"
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("D38").Value = "11.04"
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("D39").Value = "14.43"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "0.09268"
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.7782"
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.469"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "16.09"
$ws.Range("E43").Value = "  +4.67%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "2.614"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.7230"
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("D46").Value = "4.208"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "139.74"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").Value = "0.07989"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "89.81"
$ws.Range("E51").Value = "  +2.20%  "
